$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.611.09"
$ws.Range("E2").Value = "  -0.39%  "

# Row 3
$ws.Range("D3").Value = "3.261.60"
$ws.Range("E3").Value = "  +4.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'595.82"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("D6").Value = "'141.65"
$ws.Range("E6").Value = "  +1.68%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "3.255.87"
$ws.Range("E8").Value = "  +4.02%  "

# Row 9
$ws.Range("E9").Value = "  -1.45%  "

# Row 10
$ws.Range("E10").Value = "  +0.04%  "

# Row 11
$ws.Range("E11").Value = "  +1.38%  "

# Row 12
$ws.Range("E12").Value = "  +1.16%  "

# Row 14
$ws.Range("D14").Value = "'34.47"
$ws.Range("E14").Value = "  +0.12%  "

# Row 15
$ws.Range("D15").Value = "3.791.49"
$ws.Range("E15").Value = "  +3.95%  "

# Row 16
$ws.Range("E16").Value = "  -0.16%  "

# Row 17
$ws.Range("D17").Value = "3.255.30"
$ws.Range("E17").Value = "  +3.88%  "

# Row 18
$ws.Range("D18").Value = "63.618.16"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19
$ws.Range("D19").Value = "'6.81"
$ws.Range("E19").Value = "  +0.76%  "

# Row 20
$ws.Range("D20").Value = "'478.71"
$ws.Range("E20").Value = "  -0.85%  "

# Row 21
$ws.Range("D21").Value = "'14.28"
$ws.Range("E21").Value = "  -1.48%  "

# Row 22
$ws.Range("D22").Value = "'0.733"
$ws.Range("E22").Value = "  +4.23%  "

# Row 23
$ws.Range("D23").Value = "'7.99"
$ws.Range("E23").Value = "  +3.78%  "

# Row 24
$ws.Range("D24").Value = "'83.99"
$ws.Range("E24").Value = "  -4.53%  "

# Row 25
$ws.Range("E25").Value = "  +2.08%  "

# Row 26
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("D27").Value = "'2.75"
$ws.Range("E27").Value = "  +0.22%  "

# Row 28
$ws.Range("D28").Value = "'7.18"
$ws.Range("E28").Value = "  +3.64%  "

# Row 29
$ws.Range("E29").Value = "  -0.44%  "

# Row 30
$ws.Range("E30").Value = "  +4.67%  "

# Row 31
$ws.Range("D31").Value = "'27.83"
$ws.Range("E31").Value = "  +2.59%  "

# Row 32
$ws.Range("E32").Value = "  +0.03%  "

# Row 33
$ws.Range("E33").Value = "  -3.82%  "

# Row 34
$ws.Range("E34").Value = "  -1.06%  "

# Row 35
$ws.Range("E35").Value = "  -0.37%  "

# Row 36
$ws.Range("D36").Value = "'5.94"
$ws.Range("E36").Value = "  -1.19%  "

# Row 37
$ws.Range("D37").Value = "'52.80"
$ws.Range("E37").Value = "  +0.42%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0721"
$ws.Range("E38").Value = "  -2.14%  "

# Row 40
$ws.Range("D40").Value = "'422.60"
$ws.Range("E40").Value = "  -1.17%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'2.77"
$ws.Range("E41").Value = "  -2.64%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.992.60"
$ws.Range("E42").Value = "  +4.14%  "

# Row 43
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.41"
$ws.Range("E43").Value = "  +1.22%  "

# Row 44
$ws.Range("E44").Value = "  -8.19%  "

# Row 45
$ws.Range("E45").Value = "  +2.29%  "

# Row 46
$ws.Range("E46").Value = "  +1.66%  "

# Row 48
$ws.Range("E48").Value = "  -1.31%  "

# Row 49
$ws.Range("D49").Value = "'26.00"
$ws.Range("E49").Value = "  +1.79%  "

# Row 50
$ws.Range("D50").Value = "'0.115"
$ws.Range("E50").Value = "  +0.43%  "

# Row 51
$ws.Range("D51").Value = "'122.71"
$ws.Range("E51").Value = "  +1.91%  "
